{"js": "const replacements = [\n  [\"646\u00f73=215, 1\", \"604\u00f73=201, 1\"],\n  [\"554\u00f73=184, 2\", \"961\u00f75=192, 1\"],\n  [\"618\u00f73=206, 0\", \"951\u00f79=105, 6\"],\n  [\"539\u00f75=107, 4\", \"944\u00f72=472, 0\"],\n  [\"935\u00f79=103, 8\", \"761\u00f79=84, 5\"],\n  [\"640\u00f79=71, 1\", \"976\u00f78=122, 0\"],\n  [\"281\u00f75=56, 1\", \"309\u00f78=38, 5\"],\n  [\"318\u00f72=159, 0\", \"610\u00f78=76, 2\"],\n  [\"654\u00f79=72, 6\", \"209\u00f77=29, 6\"],\n  [\"289\u00f79=32, 1\", \"831\u00f76=138, 3\"],\n  [\"134\u00f77=19, 1\", \"561\u00f76=93, 3\"],\n  [\"613\u00f75=122, 3\", \"298\u00f74=74, 2\"],\n  [\"638\u00f79=70, 8\", \"359\u00f72=179, 1\"],\n  [\"221\u00f78=27, 5\", \"901\u00f72=450, 1\"],\n  [\"598\u00f79=66, 4\", \"962\u00f73=320, 2\"],\n  [\"788\u00f77=112, 4\", \"142\u00f72=71, 0\"],\n  [\"882\u00f72=441, 0\", \"373\u00f79=41, 4\"],\n  [\"587\u00f77=83, 6\", \"770\u00f74=192, 2\"],\n  [\"808\u00f73=269, 1\", \"326\u00f73=108, 2\"],\n  [\"238\u00f76=39, 4\", \"734\u00f78=91, 6\"],\n  [\"761\u00f73=253, 2\", \"228\u00f78=28, 4\"],\n  [\"699\u00f76=116, 3\", \"810\u00f79=90, 0\"],\n  [\"512\u00f74=128, 0\", \"989\u00f77=141, 2\"],\n  [\"683\u00f72=341, 1\", \"840\u00f77=120, 0\"],\n  [\"972\u00f74=243, 0\", \"881\u00f74=220, 1\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length !== 1) {\n    throw new Error('Expected exactly 1 match for \"' + oldText + '\", found ' + results.items.length);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"646\u00f73=215, 1\", \"604\u00f73=201, 1\"),\n    @(\"554\u00f73=184, 2\", \"961\u00f75=192, 1\"),\n    @(\"618\u00f73=206, 0\", \"951\u00f79=105, 6\"),\n    @(\"539\u00f75=107, 4\", \"944\u00f72=472, 0\"),\n    @(\"935\u00f79=103, 8\", \"761\u00f79=84, 5\"),\n    @(\"640\u00f79=71, 1\", \"976\u00f78=122, 0\"),\n    @(\"281\u00f75=56, 1\", \"309\u00f78=38, 5\"),\n    @(\"318\u00f72=159, 0\", \"610\u00f78=76, 2\"),\n    @(\"654\u00f79=72, 6\", \"209\u00f77=29, 6\"),\n    @(\"289\u00f79=32, 1\", \"831\u00f76=138, 3\"),\n    @(\"134\u00f77=19, 1\", \"561\u00f76=93, 3\"),\n    @(\"613\u00f75=122, 3\", \"298\u00f74=74, 2\"),\n    @(\"638\u00f79=70, 8\", \"359\u00f72=179, 1\"),\n    @(\"221\u00f78=27, 5\", \"901\u00f72=450, 1\"),\n    @(\"598\u00f79=66, 4\", \"962\u00f73=320, 2\"),\n    @(\"788\u00f77=112, 4\", \"142\u00f72=71, 0\"),\n    @(\"882\u00f72=441, 0\", \"373\u00f79=41, 4\"),\n    @(\"587\u00f77=83, 6\", \"770\u00f74=192, 2\"),\n    @(\"808\u00f73=269, 1\", \"326\u00f73=108, 2\"),\n    @(\"238\u00f76=39, 4\", \"734\u00f78=91, 6\"),\n    @(\"761\u00f73=253, 2\", \"228\u00f78=28, 4\"),\n    @(\"699\u00f76=116, 3\", \"810\u00f79=90, 0\"),\n    @(\"512\u00f74=128, 0\", \"989\u00f77=141, 2\"),\n    @(\"683\u00f72=341, 1\", \"840\u00f77=120, 0\"),\n    @(\"972\u00f74=243, 0\", \"881\u00f74=220, 1\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $ok = $find.Execute([ref]$oldText, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$newText, 2)\n    if (-not $ok) {\n        throw \"Replacement failed for: $oldText\"\n    }\n}"}
